$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B ("Count Raw Reads (Allyson)") was re-derived from the
# cutadapt-trimmed read counts already present in column D
# ("Count Raw Reads (Michael)") for the test file used to check
# against MS's counts. Copy D2:D97 -> B2:B97.
for ($r = 2; $r -le 97; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 4).Value2
}

# Restore the selection left active in the sheet after the edit.
$ws.Range("A2:B98").Select()
